# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Thu Oct 24 22:30:09 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.377.26'
$ws.Range('E2').Value = '  +3.01%  '

# Row 3
$ws.Range('D3').Value = '2.537.80'
$ws.Range('E3').Value = '  +1.34%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.05'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.12'
$ws.Range('E6').Value = '  +4.59%  '

# Row 7
$ws.Range('E7').Value = '  -0.07%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.533'
$ws.Range('E8').Value = '  +1.70%  '

# Row 9
$ws.Range('D9').Value = '2.537.82'
$ws.Range('E9').Value = '  +1.35%  '

# Row 10
$ws.Range('E10').Value = '  +2.67%  '

# Row 11
$ws.Range('E11').Value = '  +2.91%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.14'
$ws.Range('E12').Value = '  +0.61%  '

# Row 13
$ws.Range('E13').Value = '  +0.24%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.00'
$ws.Range('E14').Value = '  +1.39%  '

# Row 15
$ws.Range('D15').Value = '3.002.47'
$ws.Range('E15').Value = '  +1.53%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000180'
$ws.Range('E16').Value = '  +2.53%  '

# Row 17
$ws.Range('D17').Value = '68.282.54'
$ws.Range('E17').Value = '  +3.10%  '

# Row 18
$ws.Range('D18').Value = '2.551.22'
$ws.Range('E18').Value = '  +1.86%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.04'
$ws.Range('E19').Value = '  +4.60%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.56'
$ws.Range('E20').Value = '  +2.98%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '369.46'
$ws.Range('E21').Value = '  +6.49%  '

# Row 22
$ws.Range('E22').Value = '  +0.91%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.74'
$ws.Range('E23').Value = '  +2.59%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.94'
$ws.Range('E24').Value = '  -0.70%  '

# Row 25
$ws.Range('E25').Value = '  -0.01%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.80'
$ws.Range('E26').Value = '  +1.62%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.26'
$ws.Range('E27').Value = '  +3.42%  '

# Row 28
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.30%  '

# Row 29
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '2.661.50'
$ws.Range('E29').Value = '  +1.17%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0999'
$ws.Range('E30').Value = '  +2.87%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '543.20'
$ws.Range('E31').Value = '  +3.70%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.32'
$ws.Range('E32').Value = '  +3.07%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.34'
$ws.Range('E33').Value = '  +2.08%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  +2.76%  '

# Row 35
$ws.Range('E35').Value = '  -0.25%  '

# Row 36
$ws.Range('E36').Value = '  +0.00%  '

# Row 37
$ws.Range('E37').Value = '  +0.64%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '157.48'
$ws.Range('E38').Value = '  +0.58%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.89'
$ws.Range('E39').Value = '  +1.76%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.70'
$ws.Range('E40').Value = '  +1.79%  '

# Row 41
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.357'
$ws.Range('E41').Value = '  +0.82%  '

# Row 42
$ws.Range('E42').Value = '  +2.13%  '

# Row 43
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.24'
$ws.Range('E43').Value = '  +3.38%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.58'
$ws.Range('E44').Value = '  +3.79%  '

# Row 45
$ws.Range('E45').Value = '  -0.01%  '

# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '148.01'
$ws.Range('E46').Value = '  -0.17%  '

# Row 47
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.564'
$ws.Range('E47').Value = '  +1.44%  '

# Row 48
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0280'
$ws.Range('E48').Value = '  +3.41%  '

# Row 49
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.74'
$ws.Range('E49').Value = '  +1.91%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.71'
$ws.Range('E50').Value = '  +0.23%  '

# Row 51
$ws.Range('E51').Value = '  +1.21%  '
